# The deck's single "live" theme (ppt/theme/theme2.xml, the one actually
# referenced by the slide master / slides / presentation.xml) currently
# carries the "Integral" color scheme. The target edit swaps it for the
# stock "Office Theme" color scheme (the colors that used to live in
# ppt/theme/theme1.xml, which only the notes master pointed at).
#
# PowerPoint's COM object model exposes the twelve theme colors through
# ThemeColorScheme (indices 1-12 = Dark1, Light1, Dark2, Light2,
# Accent1..Accent6, Hyperlink, FollowedHyperlink). Driving that object
# from any slide updates the shared theme part used by the whole deck.

$p = $ppt.ActivePresentation
$cs = $p.Slides.Item(1).ThemeColorScheme

$cs.Item(1).RGB  = 0        # Dark1      -> 000000
$cs.Item(2).RGB  = 16777215 # Light1     -> FFFFFF
$cs.Item(3).RGB  = 6968388  # Dark2      -> 44546A
$cs.Item(4).RGB  = 15132391 # Light2     -> E7E6E6
$cs.Item(5).RGB  = 13998939 # Accent1    -> 5B9BD5
$cs.Item(6).RGB  = 3243501  # Accent2    -> ED7D31
$cs.Item(7).RGB  = 10855845 # Accent3    -> A5A5A5
$cs.Item(8).RGB  = 49407    # Accent4    -> FFC000
$cs.Item(9).RGB  = 12874308 # Accent5    -> 4472C4
$cs.Item(10).RGB = 4697456  # Accent6    -> 70AD47
$cs.Item(11).RGB = 12673797 # Hyperlink  -> 0563C1
$cs.Item(12).RGB = 7491477  # FollowedHyperlink -> 954F72
